$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = "ambiguous"
$ws.Range("C9").Value = "equivocal"
$ws.Range("D9").Value = "adj."

$ws.Range("B10").Value = "umambiguous"
$ws.Range("C10").Value = "unequivocal"
$ws.Range("D10").Value = "adj."

$ws.Range("B11").Select()
